# Fix the typo in the "Browser:" label and update the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct cell content change: "Browser:" -> "Broswer:" (typo introduced by the author)
$ws.Range("A2").Value = "Broswer:"

# Update the current selection from H20 to A2
$ws.Range("A2").Select()
